$d = $word.ActiveDocument

# 1) Merge the two runs around the removed "_GoBack" bookmark into a single
#    continuous sentence (the bookmark that sat between them disappears once
#    the text becomes one contiguous run again).
$d.Content.Find.Execute(
    "搭建服务单、保单、理赔等领域服务，利用其上下文的解耦和内聚进一步解决理赔补偿、组件开放、保险抽象等难题",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "搭建服务单、保单、理赔等领域服务，利用其上下文的解耦和内聚进一步解决理赔补偿、组件开放、保险抽象等难题",
    2) | Out-Null

# 2) Reword the logistics value-add intro sentence.
$d.Content.Find.Execute(
    " 物流增值业务主要包括面向商家的有赞寄件和面向消费者的上门取件，其核心都是通过与三方物流间的差价赚取利润。",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " 物流增值业务包括有赞寄件（面向商家）和上门取件（面向消费者），其核心都是通过为提供更便利更具性价比的服务，来赚取与三方物流间的差价利润。",
    2) | Out-Null

# 3) Expand bullet 1 about 有赞寄件.
$d.Content.Find.Execute(
    "  1. 有赞寄件是在交易下单后物流发货、结算的核心，其通过物流商运营、运费定价等构造整个物流计费结算模型",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "  1. 有赞寄件是交易下单后物流发货、结算的核心，其通过物流商运营、运费定价等构造整个物流计费结算模型，并通过限流、Bond分布式锁及脚本对账等方式保证最终一致性",
    2) | Out-Null

# 4) Reword bullet 2 about 上门取件.
$d.Content.Find.Execute(
    "2. 上门取件是作为交易逆向售后的核心服务，其通过状态机+最终一致保证取件单、三方物流单、交易单之间状态流转",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "2. 上门取件是交易逆向售后的核心服务，状态机 + Seata 驱动取件单、三方物流单、交易单、包运费补贴之间状态一致",
    2) | Out-Null

# 5) Re-insert the "_GoBack" bookmark right after bullet 3's sentence
#    (before its trailing line-break run).
$findRange = $d.Content
$findRange.Find.Execute(
    "3. 搭建了与三方交互的物流基础，并通过心跳、监控、预警等方式维系其稳定",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bookmarkRange = $d.Range($findRange.End, $findRange.End)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null
